$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B. This shifts the existing "Valor"
# column (and its data) from B to C, freeing up column B for the new
# "Variável" column.
$ws.Columns("B:B").Insert()

# --- Header row ---
$ws.Range("B1").Value = "Variável"
$ws.Range("C1").Value = "Valor"
$ws.Range("D1").Value = "Colocação"

# Match the header style (bold/centered) used by A1/C1.
$ws.Range("A1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Column B: "Variável" label for every data row ---
$ws.Range("B2:B10").Value = "Diferença 2021-2012"

# --- Column D: "Colocação" ranking, only for rows 2-8 ---
$ws.Range("D2").Value = "1º"
$ws.Range("D3").Value = "2º"
$ws.Range("D4").Value = "3º"
$ws.Range("D5").Value = "4º"
$ws.Range("D6").Value = "5º"
$ws.Range("D7").Value = "6º"
$ws.Range("D8").Value = "24º"
